# "Update countries & provincias Spain"
#
# Daily refresh of the COVID "Pais" sheet:
#   - the "Datos actualizados a ..." timestamp caption moves forward,
#   - most country rows receive updated totals,
#   - a handful of countries are re-ranked by total cases, which (since the
#     table is sorted descending by column B) swaps some adjacent rows'
#     labels along with their figures:
#       * Peru      overtakes España               (rows 8 / 9)
#       * Barein    overtakes Rumania               (rows 50 / 51)
#       * Bulgaria  overtakes Gabon and Guinea      (rows 87 / 88 / 89)
#       * Zimbabue  overtakes San Marino/Malta/Togo (rows 150-153)
#       * Islas Malvinas overtakes Groenlandia      (rows 209 / 210)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Header caption (row 1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 5 de Julio de 2020 a las 00:03"

# Estados Unidos (row 4) - refreshed totals
$ws.Cells.Item(4, 2).Value = 2933428
$ws.Cells.Item(4, 3).Value = 42840
$ws.Cells.Item(4, 4).Value = 1256161
$ws.Cells.Item(4, 5).Value = 1544966
$ws.Cells.Item(4, 7).Value = 237
$ws.Cells.Item(4, 8).Value = 132301

# Brasil (row 5) - refreshed totals
$ws.Cells.Item(5, 2).Value = 1577004
$ws.Cells.Item(5, 3).Value = 33663
$ws.Cells.Item(5, 5).Value = 534124
$ws.Cells.Item(5, 7).Value = 1011
$ws.Cells.Item(5, 8).Value = 64265

# Peru now ranks above España -> row 8 becomes Peru (fresh data)
$ws.Cells.Item(8, 1).Value = "Peru"
$ws.Cells.Item(8, 2).Value = 299080
$ws.Cells.Item(8, 3).Value = 3481
$ws.Cells.Item(8, 4).Value = 189621
$ws.Cells.Item(8, 5).Value = 99047
$ws.Cells.Item(8, 7).Value = 186
$ws.Cells.Item(8, 8).Value = 10412

# ... and row 9 becomes España (carrying the old row-8 figures down)
$ws.Cells.Item(9, 1).Value = "España"
$ws.Cells.Item(9, 2).Value = 297625
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = 0
$ws.Cells.Item(9, 8).Value = 28385

# Alemania (row 18) - refreshed totals
$ws.Cells.Item(18, 2).Value = 197408
$ws.Cells.Item(18, 3).Value = 408
$ws.Cells.Item(18, 5).Value = 7027
$ws.Cells.Item(18, 7).Value = 8
$ws.Cells.Item(18, 8).Value = 9081

# Canada (row 23) - refreshed totals
$ws.Cells.Item(23, 2).Value = 105316
$ws.Cells.Item(23, 3).Value = 225
$ws.Cells.Item(23, 4).Value = 68990
$ws.Cells.Item(23, 5).Value = 27652
$ws.Cells.Item(23, 7).Value = 11
$ws.Cells.Item(23, 8).Value = 8674

# Ecuador (row 32) - refreshed totals
$ws.Cells.Item(32, 2).Value = 61535
$ws.Cells.Item(32, 3).Value = 878
$ws.Cells.Item(32, 4).Value = 28507
$ws.Cells.Item(32, 5).Value = 28259
$ws.Cells.Item(32, 7).Value = 69
$ws.Cells.Item(32, 8).Value = 4769

# Portugal (row 41) - refreshed totals
$ws.Cells.Item(41, 2).Value = 43569
$ws.Cells.Item(41, 3).Value = 413
$ws.Cells.Item(41, 5).Value = 13192

# Israel (row 49) - refreshed totals
$ws.Cells.Item(49, 2).Value = 29170
$ws.Cells.Item(49, 3).Value = 1115
$ws.Cells.Item(49, 4).Value = 17816
$ws.Cells.Item(49, 5).Value = 11024

# Barein now ranks above Rumania -> row 50 becomes Barein (fresh data)
$ws.Cells.Item(50, 1).Value = "Barein"
$ws.Cells.Item(50, 2).Value = 28857
$ws.Cells.Item(50, 3).Value = 447
$ws.Cells.Item(50, 4).Value = 23959
$ws.Cells.Item(50, 5).Value = 4802
$ws.Cells.Item(50, 7).Value = 1
$ws.Cells.Item(50, 8).Value = 96

# ... and row 51 becomes Rumania (carrying the old row-50 figures down)
$ws.Cells.Item(51, 1).Value = "Rumania"
$ws.Cells.Item(51, 2).Value = 28582
$ws.Cells.Item(51, 3).Value = 416
$ws.Cells.Item(51, 4).Value = 19854
$ws.Cells.Item(51, 5).Value = 6997
$ws.Cells.Item(51, 7).Value = 23
$ws.Cells.Item(51, 8).Value = 1731

# Chequia (row 69) - refreshed totals
$ws.Cells.Item(69, 5).Value = 4169
$ws.Cells.Item(69, 8).Value = 352

# Bulgaria now ranks above Gabon -> row 87 becomes Bulgaria (fresh data)
$ws.Cells.Item(87, 1).Value = "Bulgaria"
$ws.Cells.Item(87, 2).Value = 5677
$ws.Cells.Item(87, 3).Value = 180
$ws.Cells.Item(87, 4).Value = 2898
$ws.Cells.Item(87, 5).Value = 2538
$ws.Cells.Item(87, 7).Value = 2
$ws.Cells.Item(87, 8).Value = 241

# ... row 88 becomes Gabon (carrying the old row-87 figures down)
$ws.Cells.Item(88, 1).Value = "Gabon"
$ws.Cells.Item(88, 2).Value = 5620
$ws.Cells.Item(88, 4).Value = 2555
$ws.Cells.Item(88, 5).Value = 3021
$ws.Cells.Item(88, 8).Value = 44

# ... row 89 becomes Guinea (carrying the old row-88 figures down)
$ws.Cells.Item(89, 1).Value = "Guinea"
$ws.Cells.Item(89, 2).Value = 5521
$ws.Cells.Item(89, 4).Value = 4446
$ws.Cells.Item(89, 5).Value = 1042
$ws.Cells.Item(89, 8).Value = 33

# Estado de Palestina (row 98) - refreshed totals
$ws.Cells.Item(98, 5).Value = 3359
$ws.Cells.Item(98, 7).Value = 2
$ws.Cells.Item(98, 8).Value = 13

# Zimbabue now ranks above San Marino -> row 150 becomes Zimbabue (fresh data)
$ws.Cells.Item(150, 1).Value = "Zimbabue"
$ws.Cells.Item(150, 3).Value = 73
$ws.Cells.Item(150, 4).Value = 181
$ws.Cells.Item(150, 5).Value = 509
$ws.Cells.Item(150, 7).Value = 1
$ws.Cells.Item(150, 8).Value = 8

# ... row 151 becomes San Marino (carrying the old row-150 figures down)
$ws.Cells.Item(151, 1).Value = "San Marino"
$ws.Cells.Item(151, 2).Value = 698
$ws.Cells.Item(151, 4).Value = 656
$ws.Cells.Item(151, 5).Value = 0
$ws.Cells.Item(151, 8).Value = 42

# ... row 152 becomes Malta (carrying the old row-151 figures down)
$ws.Cells.Item(152, 1).Value = "Malta"
$ws.Cells.Item(152, 2).Value = 672
$ws.Cells.Item(152, 4).Value = 650
$ws.Cells.Item(152, 5).Value = 13
$ws.Cells.Item(152, 8).Value = 9

# ... row 153 becomes Togo (carrying the old row-152 figures down)
$ws.Cells.Item(153, 1).Value = "Togo"
$ws.Cells.Item(153, 2).Value = 671
$ws.Cells.Item(153, 4).Value = 424
$ws.Cells.Item(153, 5).Value = 233
$ws.Cells.Item(153, 8).Value = 14

# Gambia (row 190) - refreshed totals
$ws.Cells.Item(190, 2).Value = 57
$ws.Cells.Item(190, 3).Value = 2
$ws.Cells.Item(190, 5).Value = 28

# San Cristobal y Nieves (row 208) - refreshed totals
$ws.Cells.Item(208, 2).Value = 16
$ws.Cells.Item(208, 3).Value = 1
$ws.Cells.Item(208, 5).Value = 1

# Islas Malvinas now ranks above Groenlandia (same totals, labels only)
$ws.Cells.Item(209, 1).Value = "Islas Malvinas"
$ws.Cells.Item(210, 1).Value = "Groenlandia"
